# Ultimate_Job_Search_Tracker.xlsx edit script
# Reproduces: new "Compensation" column on "Job Applications", refreshed
# ATS-focused tip text on "Tips & Guidance", plus the various column-width /
# row-height / selection / active-sheet cosmetic changes captured in the diff.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: Job Applications
# Insert a new "Compensation" column between "Contact Person" and "Notes".
# ---------------------------------------------------------------------
$wsJobs = $wb.Worksheets.Item(1)

# Shift the existing "Notes" header (I1) into the new last column (J1),
# carrying its style along, then overwrite I1 with the new header text.
$wsJobs.Range("I1").Copy($wsJobs.Range("J1"))
$wsJobs.Range("I1").Value = "Compensation"

# Taller header row to fit the wrapped text of the extra column.
$wsJobs.Rows.Item(1).RowHeight = 36.6

# New column widths for H:J; this also replaces the old single-column (I)
# "JD/Email Link" width definition.
$wsJobs.Columns.Item(8).ColumnWidth = 18.498697916666668
$wsJobs.Columns.Item(9).ColumnWidth = 17.498697916666668
$wsJobs.Columns.Item(10).ColumnWidth = 31.498697916666668

$wsJobs.Range("J2").Select()

# ---------------------------------------------------------------------
# Sheet 2: Preparation Plan
# Give column F (already holding the "Notes" header) an explicit width.
# ---------------------------------------------------------------------
$wsPrep = $wb.Worksheets.Item(2)
$wsPrep.Columns.Item(6).ColumnWidth = 26.166666666666668

# ---------------------------------------------------------------------
# Sheet 3: Daily To-Do
# ---------------------------------------------------------------------
$wsDaily = $wb.Worksheets.Item(3)
$wsDaily.Range("F2").Select()

# ---------------------------------------------------------------------
# Sheet 4: Referrals
# ---------------------------------------------------------------------
$wsReferrals = $wb.Worksheets.Item(4)
$wsReferrals.Range("G3").Select()

# ---------------------------------------------------------------------
# Sheet 5: Mental Health & Routine
# ---------------------------------------------------------------------
$wsMental = $wb.Worksheets.Item(5)
$wsMental.Range("E28").Select()

# ---------------------------------------------------------------------
# Sheet 6: Tips & Guidance
# Refresh three tip strings to the new ATS-focused wording and resize
# the category/tip columns.
# ---------------------------------------------------------------------
$wsTips = $wb.Worksheets.Item(6)
$wsTips.Range("B4").Value = "Set a 4-week roadmap: ATS friendly resume, profile update, applications, referrals, interview prep."
$wsTips.Range("B5").Value = "Make it ATS friendly, focus on outcomes. Use metrics. Avoid generic summaries. "
$wsTips.Range("B9").Value = "Apply to 10 jobs, follow up 2 referrals, practice 5 interview question. Try to avoid easy apply as there would be already a lot of applications"

$wsTips.Columns.Item(1).ColumnWidth = 25.498697916666668
$wsTips.Columns.Item(2).ColumnWidth = 113.94401041666667

$wsTips.Range("B10").Select()

# ---------------------------------------------------------------------
# Sheet 7: Networking & Events
# ---------------------------------------------------------------------
$wsNetworking = $wb.Worksheets.Item(7)
$wsNetworking.Range("E2:E13").Select()

# ---------------------------------------------------------------------
# Sheet 8: Skills Tracker
# Becomes the active / selected sheet.
# ---------------------------------------------------------------------
$wsSkills = $wb.Worksheets.Item(8)
$wsSkills.Activate()
$wsSkills.Range("E4").Select()
